$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename headers to be R-friendly (lowercase, underscores instead of spaces)
$ws.Range("A1").Value = "date"
$ws.Range("B1").Value = "total_calories_burned"
$ws.Range("C1").Value = "daily_step_count"
